# Update fuel_summary sheet with latest P&L calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dataset: flight_id, fuel_burn_rate_gph, flight_hours, fuel_price_per_gal, fuel_cost
$data = @(
    @("FL6679", 850,  2.03, 3.409, 5882.23),
    @("FL1976", 850,  3.68, 2.427, 7591.66),
    @("FL4876", 850,  4.95, 2.427, 10211.6),
    @("FL1198", 850,  2.03, 2.427, 4187.79),
    @("FL3066", 850,  2.5,  2.427, 5157.38),
    @("FL5388", 850,  3.68, 2.427, 7591.66),
    @("FL5301", 1000, 2.5,  2.427, 6067.5),
    @("FL7875", 1000, 4.95, 2.427, 12013.65),
    @("FL6369", 1000, 3.68, 2.427, 8931.360000000001),
    @("FL7466", 1000, 3.68, 2.427, 8931.360000000001),
    @("FL3108", 1000, 2.5,  2.427, 6067.5),
    @("FL4736", 1000, 1.45, 2.427, 3519.15),
    @("FL5961", 1350, 7.25, 2.427, 23754.26),
    @("FL2734", 1350, 7.25, 2.427, 23754.26),
    @("FL1013", 1350, 7.25, 2.427, 23754.26),
    @("FL5705", 890,  2.03, 3.409, 6159.04),
    @("FL1524", 890,  3.68, 2.427, 7948.91),
    @("FL5421", 890,  4.95, 2.427, 10692.15),
    @("FL6133", 890,  4.95, 2.427, 10692.15),
    @("FL8299", 500,  1.45, 2.427, 1759.58),
    @("FL5785", 790,  1.45, 3.409, 3905.01),
    @("FL5235", 790,  4.95, 2.427, 9490.780000000001),
    @("FL5689", 790,  2.03, 2.427, 3892.18),
    @("FL8121", 790,  1.45, 2.427, 2780.13),
    @("FL1990", 1350, 7.25, 2.427, 23754.26)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
